$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1) to snake_case English names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case the lowercase Spanish connector words ("de", "del", "la", "las", "el", "los", "y")
# in municipality/state names throughout column A and B
$ws.Range("B6").Value = "Pabellón De Arteaga"
$ws.Range("B7").Value = "Rincón De Romos"
$ws.Range("B8").Value = "San José De Gracia"
$ws.Range("B21").Value = "Amatenango De La Frontera"
$ws.Range("B26").Value = "Chiapa De Corzo"
$ws.Range("B29").Value = "Comitán De Domínguez"
$ws.Range("B41").Value = "Mazapa De Madero"
$ws.Range("B42").Value = "Montecristo De Guerrero"
$ws.Range("B48").Value = "San Cristóbal De Las Casas"
$ws.Range("B72").Value = "Hidalgo Del Parral"
$ws.Range("B80").Value = "San Francisco De Borja"
$ws.Range("B81").Value = "San Francisco Del Oro"
$ws.Range("B85").Value = "Valle De Zaragoza"
$ws.Range("B97").Value = "San Juan De Sabinas"
$ws.Range("B107").Value = "Villa De Álvarez"
$ws.Range("A109").Value = "Ciudad De México"
$ws.Range("B113").Value = "Cuajimalpa De Morelos"
$ws.Range("B127").Value = "Coneto De Comonfort"
$ws.Range("B140").Value = "Nombre De Dios"
$ws.Range("B149").Value = "San Juan Del Río"
$ws.Range("A157").Value = "Estado De México"
$ws.Range("B157").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B160").Value = "Almoloya De Alquisiras"
$ws.Range("B161").Value = "Almoloya De Juárez"
$ws.Range("B167").Value = "Atizapán De Zaragoza"
$ws.Range("B171").Value = "Chapa De Mota"
$ws.Range("B173").Value = "Coacalco De Berriozábal"
$ws.Range("B178").Value = "Ecatepec De Morelos"
$ws.Range("B182").Value = "Ixtapan De La Sal"
$ws.Range("B191").Value = "Naucalpan De Juárez"
$ws.Range("B197").Value = "San Felipe Del Progreso"
$ws.Range("B198").Value = "San Martín De Las Pirámides"
$ws.Range("B212").Value = "Tlalnepantla De Baz"
$ws.Range("B217").Value = "Valle De Bravo"
$ws.Range("B218").Value = "Villa De Allende"
$ws.Range("B227").Value = "San Miguel De Allende"
$ws.Range("B228").Value = "Apaseo El Alto"
$ws.Range("B229").Value = "Apaseo El Grande"
$ws.Range("B235").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B239").Value = "Jaral Del Progreso"
$ws.Range("B246").Value = "Purísima Del Rincón"
$ws.Range("B250").Value = "San Diego De La Unión"
$ws.Range("B252").Value = "San Francisco Del Rincón"
$ws.Range("B254").Value = "San Luis De La Paz"
$ws.Range("B256").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B258").Value = "Silao De La Victoria"
$ws.Range("B262").Value = "Valle De Santiago"
$ws.Range("B268").Value = "Acapulco De Juárez"
$ws.Range("B270").Value = "Ajuchitlán Del Progreso"
$ws.Range("B274").Value = "Atoyac De Álvarez"
$ws.Range("B275").Value = "Ayutla De Los Libres"
$ws.Range("B277").Value = "Buenavista De Cuéllar"
$ws.Range("B278").Value = "Chilapa De Álvarez"
$ws.Range("B279").Value = "Chilpancingo De Los Bravo"
$ws.Range("B280").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B284").Value = "Coyuca De Benítez"
$ws.Range("B285").Value = "Coyuca De Catalán"
$ws.Range("B288").Value = "Cuetzala Del Progreso"
$ws.Range("B289").Value = "Cutzamala De Pinzón"
$ws.Range("B294").Value = "Huitzuco De Los Figueroa"
$ws.Range("B295").Value = "Iguala De La Independencia"
$ws.Range("B296").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B297").Value = "Zihuatanejo De Azueta"
$ws.Range("B311").Value = "Taxco De Alarcón"
$ws.Range("B313").Value = "Técpan De Galeana"
$ws.Range("B315").Value = "Tepecoacuilco De Trujano"
$ws.Range("B317").Value = "Tixtla De Guerrero"
$ws.Range("B320").Value = "Tlapa De Comonfort"
$ws.Range("B332").Value = "Atotonilco El Grande"
$ws.Range("B335").Value = "Cuautepec De Hinojosa"
$ws.Range("B339").Value = "Huejutla De Reyes"
$ws.Range("B342").Value = "Jacala De Ledezma"
$ws.Range("B347").Value = "Pachuca De Soto"
$ws.Range("B350").Value = "Santiago De Anaya"
$ws.Range("B352").Value = "Tenango De Doria"
$ws.Range("B354").Value = "Tepehuacán De Guerrero"
$ws.Range("B355").Value = "Tezontepec De Aldama"
$ws.Range("B360").Value = "Tula De Allende"
$ws.Range("B361").Value = "Tulancingo De Bravo"
$ws.Range("B363").Value = "Zacualtipán De Ángeles"
$ws.Range("B364").Value = "Zapotlán De Juárez"
$ws.Range("B368").Value = "Ahualulco De Mercado"
$ws.Range("B372").Value = "Atemajac De Brizuela"
$ws.Range("B373").Value = "Atotonilco El Alto"
$ws.Range("B374").Value = "Autlán De Navarro"
$ws.Range("B389").Value = "Encarnación De Díaz"
$ws.Range("B392").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B393").Value = "Ixtlahuacán Del Río"
$ws.Range("B399").Value = "Lagos De Moreno"
$ws.Range("B409").Value = "San Cristóbal De La Barranca"
$ws.Range("B410").Value = "San Diego De Alejandría"
$ws.Range("B412").Value = "San Juan De Los Lagos"
$ws.Range("B414").Value = "San Miguel El Alto"
$ws.Range("B415").Value = "Santa María De Los Ángeles"
$ws.Range("B418").Value = "Tamazula De Gordiano"
$ws.Range("B422").Value = "Tepatitlán De Morelos"
$ws.Range("B427").Value = "Unión De San Antonio"
$ws.Range("B428").Value = "Unión De Tula"
$ws.Range("B432").Value = "Yahualica De González Gallo"
$ws.Range("B433").Value = "Zacoalco De Torres"
$ws.Range("B436").Value = "Zapotlán Del Rey"
$ws.Range("B437").Value = "Zapotlán El Grande"
$ws.Range("B528").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B531").Value = "Puente De Ixtla"
$ws.Range("B535").Value = "Tetela Del Volcán"
$ws.Range("B544").Value = "Amatlán De Cañas"
$ws.Range("B547").Value = "Ixtlán Del Río"
$ws.Range("B551").Value = "Santa María Del Oro"
$ws.Range("B567").Value = "San Nicolás De Los Garza"
$ws.Range("B571").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B575").Value = "Ayoquezco De Aldama"
$ws.Range("B577").Value = "Chalcatongo De Hidalgo"
$ws.Range("B580").Value = "Coicoyán De Las Flores"
$ws.Range("B582").Value = "El Barrio De La Soledad"
$ws.Range("B583").Value = "Fresnillo De Trujano"
$ws.Range("B584").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B585").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B586").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B588").Value = "Ixtlán De Juárez"
$ws.Range("B589").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B595").Value = "Mariscala De Juárez"
$ws.Range("B596").Value = "Mártires De Tacubaya"
$ws.Range("B599").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B600").Value = "Nejapa De Madero"
$ws.Range("B601").Value = "Oaxaca De Juárez"
$ws.Range("B602").Value = "Ocotlán De Morelos"
$ws.Range("B604").Value = "Putla Villa De Guerrero"
$ws.Range("B605").Value = "Reforma De Pineda"
$ws.Range("B614").Value = "San Baltazar Yatzachi El Bajo"
$ws.Range("B619").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B660").Value = "San Pedro El Alto"
$ws.Range("B671").Value = "San Pedro Y San Pablo Ayutla"
$ws.Range("B672").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B685").Value = "Santa Cruz De Bravo"
$ws.Range("B687").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B689").Value = "Santa Inés Del Monte"
$ws.Range("B723").Value = "Santo Domingo De Morelos"
$ws.Range("B730").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B731").Value = "Tanetze De Zaragoza"
$ws.Range("B732").Value = "Tataltepec De Valdés"
$ws.Range("B733").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B734").Value = "Tlacolula De Matamoros"
$ws.Range("B735").Value = "Villa De Tututepec"
$ws.Range("B736").Value = "Villa De Zaachila"
$ws.Range("B738").Value = "Villa Sola De Vega"
$ws.Range("B739").Value = "Villa Talea De Castro"
$ws.Range("B740").Value = "Yutanduchi De Guerrero"
$ws.Range("B741").Value = "Zapotitlán Del Río"
$ws.Range("B744").Value = "Zimatlán De Álvarez"
$ws.Range("B758").Value = "Ayotoxco De Guerrero"
$ws.Range("B760").Value = "Chalchicomula De Sesma"
$ws.Range("B766").Value = "Chila De La Sal"
$ws.Range("B778").Value = "Huehuetlán El Chico"
$ws.Range("B779").Value = "Huehuetlán El Grande"
$ws.Range("B785").Value = "Izúcar De Matamoros"
$ws.Range("B796").Value = "Palmar De Bravo"
$ws.Range("B813").Value = "San Salvador El Seco"
$ws.Range("B818").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B823").Value = "Tepexi De Rodríguez"
$ws.Range("B824").Value = "Tetela De Ocampo"
$ws.Range("B827").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B847").Value = "Amealco De Bonfil"
$ws.Range("B852").Value = "Jalpan De Serra"
$ws.Range("B856").Value = "San Juan Del Río"
$ws.Range("B868").Value = "Ciudad Del Maíz"
$ws.Range("B879").Value = "Santa María Del Río"
$ws.Range("B882").Value = "Villa De Ramos"
$ws.Range("B883").Value = "Villa De Reyes"
$ws.Range("B934").Value = "Soto La Marina"
$ws.Range("B946").Value = "Contla De Juan Cuamatzi"
$ws.Range("B952").Value = "Papalotla De Xicohténcatl"
$ws.Range("B953").Value = "Sanctórum De Lázaro Cárdenas"
$ws.Range("B966").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B974").Value = "Boca Del Río"
$ws.Range("B976").Value = "Camarón De Tejeda"
$ws.Range("B990").Value = "Cosamaloapan De Carpio"
$ws.Range("B991").Value = "Cosautlán De Carvajal"
$ws.Range("B1002").Value = "Hueyapan De Ocampo"
$ws.Range("B1003").Value = "Ignacio De La Llave"
$ws.Range("B1006").Value = "Ixhuatlán De Madero"
$ws.Range("B1007").Value = "Ixhuatlán Del Café"
$ws.Range("B1008").Value = "Ixhuatlán Del Sureste"
$ws.Range("B1017").Value = "Juchique De Ferrer"
$ws.Range("B1020").Value = "Landero Y Coss"
$ws.Range("B1022").Value = "Lerdo De Tejada"
$ws.Range("B1024").Value = "Martínez De La Torre"
$ws.Range("B1025").Value = "Medellín De Bravo"
$ws.Range("B1029").Value = "Mixtla De Altamirano"
$ws.Range("B1039").Value = "Ozuluama De Mascareñas"
$ws.Range("B1042").Value = "Paso De Ovejas"
$ws.Range("B1043").Value = "Paso Del Macho"
$ws.Range("B1046").Value = "Poza Rica De Hidalgo"
$ws.Range("B1052").Value = "Sayula De Alemán"
$ws.Range("B1054").Value = "Soledad De Doblado"
$ws.Range("B1078").Value = "Vega De Alatorre"
$ws.Range("B1097").Value = "Concepción Del Oro"
$ws.Range("B1099").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B1107").Value = "Jiménez Del Teul"
$ws.Range("B1116").Value = "Nochistlán De Mejía"
$ws.Range("B1117").Value = "Noria De Ángeles"
$ws.Range("B1125").Value = "Teúl De González Ortega"
$ws.Range("B1126").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1128").Value = "Trinidad García De La Cadena"
$ws.Range("B1131").Value = "Villa De Cos"

# Fix two pairs of floating-point percentage values (1-ULP precision correction)
$ws.Range("D117").Value2 = 0.009236752552260571
$ws.Range("D119").Value2 = 0.009236752552260571
$ws.Range("D941").Value2 = 0.009479824987846381
$ws.Range("D1051").Value2 = 0.009479824987846381

# Remove trailing footnote rows 1140-1144 (source/metadata notes no longer needed)
$ws.Rows("1140:1144").Delete()
